# Scheduled-runner data refresh: updates cached market-price / profit
# figures (columns H-N) on a handful of rows across the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets. Columns:
#   H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#   K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(13, 8).Value = 10000
$ws.Cells.Item(13, 10).Value = 10000
$ws.Cells.Item(13, 12).Value = 10000
$ws.Cells.Item(13, 14).Value = -10338

$ws.Cells.Item(31, 8).Value = 587.5
$ws.Cells.Item(31, 9).Value = 587.5
$ws.Cells.Item(31, 11).Value = 1762.5
$ws.Cells.Item(31, 13).Value = -1532.5

$ws.Cells.Item(39, 8).Value = 1950.6364
$ws.Cells.Item(39, 9).Value = 135.42857
$ws.Cells.Item(39, 11).Value = 406.28571
$ws.Cells.Item(39, 13).Value = -110.28571

$ws.Cells.Item(41, 8).Value = 1132.25
$ws.Cells.Item(41, 9).Value = 504.5
$ws.Cells.Item(41, 10).Value = 2178.5
$ws.Cells.Item(41, 11).Value = 504.5
$ws.Cells.Item(41, 12).Value = 2178.5
$ws.Cells.Item(41, 13).Value = -64.5
$ws.Cells.Item(41, 14).Value = -3058.5

$ws.Cells.Item(62, 8).Value = 15690617
$ws.Cells.Item(62, 9).Value = 24246772
$ws.Cells.Item(62, 10).Value = 4333
$ws.Cells.Item(62, 11).Value = 24246772
$ws.Cells.Item(62, 12).Value = 4333
$ws.Cells.Item(62, 13).Value = -24246148
$ws.Cells.Item(62, 14).Value = -5581

$ws.Cells.Item(65, 8).Value = 15690617
$ws.Cells.Item(65, 9).Value = 24246772
$ws.Cells.Item(65, 10).Value = 4333
$ws.Cells.Item(65, 11).Value = 121233860
$ws.Cells.Item(65, 12).Value = 21665
$ws.Cells.Item(65, 13).Value = -121230740
$ws.Cells.Item(65, 14).Value = -27905

$ws.Cells.Item(94, 8).Value = 10001.667
$ws.Cells.Item(94, 9).Value = 10001.667
$ws.Cells.Item(94, 11).Value = 10001.667
$ws.Cells.Item(94, 13).Value = -9550.666999999999

$ws.Cells.Item(106, 8).Value = 3047.375
$ws.Cells.Item(106, 9).Value = 2619.75
$ws.Cells.Item(106, 10).Value = 3475
$ws.Cells.Item(106, 11).Value = 2619.75
$ws.Cells.Item(106, 12).Value = 3475
$ws.Cells.Item(106, 13).Value = -1988.75
$ws.Cells.Item(106, 14).Value = -4737

$ws.Cells.Item(116, 8).Value = 5658.8823
$ws.Cells.Item(116, 9).Value = 4807.9165
$ws.Cells.Item(116, 11).Value = 4807.9165
$ws.Cells.Item(116, 13).Value = -1365.9165

$ws.Cells.Item(129, 8).Value = 1572.76
$ws.Cells.Item(129, 9).Value = 808.0714
$ws.Cells.Item(129, 11).Value = 2424.2142
$ws.Cells.Item(129, 13).Value = 2575.7858

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 11086.087
$ws.Cells.Item(61, 9).Value = 9144
$ws.Cells.Item(61, 10).Value = 14107.111
$ws.Cells.Item(61, 11).Value = 9144
$ws.Cells.Item(61, 12).Value = 14107.111
$ws.Cells.Item(61, 13).Value = -8932
$ws.Cells.Item(61, 14).Value = -14531.111

$ws.Cells.Item(63, 8).Value = 7159.8647
$ws.Cells.Item(63, 10).Value = 7311.7144
$ws.Cells.Item(63, 12).Value = 7311.7144
$ws.Cells.Item(63, 14).Value = -8683.714400000001

$ws.Cells.Item(66, 8).Value = 7159.8647
$ws.Cells.Item(66, 10).Value = 7311.7144
$ws.Cells.Item(66, 12).Value = 36558.572
$ws.Cells.Item(66, 14).Value = -43422.572

$ws.Cells.Item(97, 8).Value = 1243053.1
$ws.Cells.Item(97, 9).Value = 1861571.6
$ws.Cells.Item(97, 11).Value = 1861571.6
$ws.Cells.Item(97, 13).Value = -1861075.6

$ws.Cells.Item(110, 8).Value = 7354324
$ws.Cells.Item(110, 9).Value = 9260093
$ws.Cells.Item(110, 11).Value = 9260093
$ws.Cells.Item(110, 13).Value = -9258048

$ws.Cells.Item(122, 8).Value = 1539.381
$ws.Cells.Item(122, 9).Value = 973.8333
$ws.Cells.Item(122, 11).Value = 2921.4999
$ws.Cells.Item(122, 13).Value = -471.4998999999998

$ws.Cells.Item(136, 8).Value = 11086.087
$ws.Cells.Item(136, 9).Value = 9144
$ws.Cells.Item(136, 10).Value = 14107.111
$ws.Cells.Item(136, 11).Value = 27432
$ws.Cells.Item(136, 12).Value = 42321.333
$ws.Cells.Item(136, 13).Value = -24882
$ws.Cells.Item(136, 14).Value = -47421.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2585.6086
$ws.Cells.Item(20, 9).Value = 2672.7693
$ws.Cells.Item(20, 10).Value = 2472.3
$ws.Cells.Item(20, 11).Value = 2672.7693
$ws.Cells.Item(20, 12).Value = 2472.3
$ws.Cells.Item(20, 13).Value = -2425.7693
$ws.Cells.Item(20, 14).Value = -2966.3

$ws.Cells.Item(99, 8).Value = 2959.2727
$ws.Cells.Item(99, 10).Value = 3247.5
$ws.Cells.Item(99, 12).Value = 3247.5
$ws.Cells.Item(99, 14).Value = -6243.5

$ws.Cells.Item(103, 8).Value = 40750.75
$ws.Cells.Item(103, 9).Value = 48000
$ws.Cells.Item(103, 10).Value = 38334.332
$ws.Cells.Item(103, 11).Value = 48000
$ws.Cells.Item(103, 12).Value = 38334.332
$ws.Cells.Item(103, 13).Value = -46828
$ws.Cells.Item(103, 14).Value = -40678.332

$ws.Cells.Item(105, 8).Value = 52646070
$ws.Cells.Item(105, 9).Value = 71447200
$ws.Cells.Item(105, 11).Value = 71447200
$ws.Cells.Item(105, 13).Value = -71445453

$ws.Cells.Item(118, 8).Value = 56249.668
$ws.Cells.Item(118, 10).Value = 56249.668
$ws.Cells.Item(118, 12).Value = 56249.668
$ws.Cells.Item(118, 14).Value = -59563.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 104.8
$ws.Cells.Item(7, 9).Value = 135.23077
$ws.Cells.Item(7, 10).Value = 48.285713
$ws.Cells.Item(7, 11).Value = 135.23077
$ws.Cells.Item(7, 12).Value = 48.285713
$ws.Cells.Item(7, 13).Value = -22.23077000000001
$ws.Cells.Item(7, 14).Value = -274.285713

$ws.Cells.Item(31, 8).Value = 47624480
$ws.Cells.Item(31, 9).Value = 100002720
$ws.Cells.Item(31, 10).Value = 7894.727
$ws.Cells.Item(31, 11).Value = 100002720
$ws.Cells.Item(31, 12).Value = 7894.727
$ws.Cells.Item(31, 13).Value = -100002425
$ws.Cells.Item(31, 14).Value = -8484.726999999999

$ws.Cells.Item(34, 8).Value = 47624480
$ws.Cells.Item(34, 9).Value = 100002720
$ws.Cells.Item(34, 10).Value = 7894.727
$ws.Cells.Item(34, 11).Value = 100002720
$ws.Cells.Item(34, 12).Value = 7894.727
$ws.Cells.Item(34, 13).Value = -100002518
$ws.Cells.Item(34, 14).Value = -8298.726999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 10104282
$ws.Cells.Item(131, 10).Value = 3692.7585
$ws.Cells.Item(131, 12).Value = 11078.2755
$ws.Cells.Item(131, 14).Value = -21158.2755

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 38755.617
$ws.Cells.Item(2, 10).Value = 125113.75
$ws.Cells.Item(2, 12).Value = 125113.75
$ws.Cells.Item(2, 14).Value = -125339.75

$ws.Cells.Item(97, 8).Value = 691.0625
$ws.Cells.Item(97, 10).Value = 681.4286
$ws.Cells.Item(97, 12).Value = 681.4286
$ws.Cells.Item(97, 14).Value = -1673.4286

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3218.568
$ws.Cells.Item(22, 9).Value = 2262.2942
$ws.Cells.Item(22, 10).Value = 3820.6667
$ws.Cells.Item(22, 11).Value = 2262.2942
$ws.Cells.Item(22, 12).Value = 3820.6667
$ws.Cells.Item(22, 13).Value = -1967.2942
$ws.Cells.Item(22, 14).Value = -4410.6667

$ws.Cells.Item(27, 8).Value = 3218.568
$ws.Cells.Item(27, 9).Value = 2262.2942
$ws.Cells.Item(27, 10).Value = 3820.6667
$ws.Cells.Item(27, 11).Value = 2262.2942
$ws.Cells.Item(27, 12).Value = 3820.6667
$ws.Cells.Item(27, 13).Value = -2155.2942
$ws.Cells.Item(27, 14).Value = -4034.6667

$ws.Cells.Item(99, 8).Value = 78000
$ws.Cells.Item(99, 10).Value = 78000
$ws.Cells.Item(99, 12).Value = 78000
$ws.Cells.Item(99, 14).Value = -83990

$ws.Cells.Item(132, 8).Value = 7579.5
$ws.Cells.Item(132, 9).Value = 4810.6
$ws.Cells.Item(132, 11).Value = 14431.8
$ws.Cells.Item(132, 13).Value = -11901.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(6, 8).Value = 2234.3333
$ws.Cells.Item(6, 9).Value = 3312.5
$ws.Cells.Item(6, 10).Value = 1695.25
$ws.Cells.Item(6, 11).Value = 3312.5
$ws.Cells.Item(6, 12).Value = 1695.25
$ws.Cells.Item(6, 13).Value = -3197.5
$ws.Cells.Item(6, 14).Value = -1925.25

$ws.Cells.Item(97, 8).Value = 27572
$ws.Cells.Item(97, 10).Value = 27572
$ws.Cells.Item(97, 12).Value = 27572
$ws.Cells.Item(97, 14).Value = -29554

$ws.Cells.Item(100, 8).Value = 1519
$ws.Cells.Item(100, 10).Value = 1871.2
$ws.Cells.Item(100, 12).Value = 3742.4
$ws.Cells.Item(100, 14).Value = -4824.4

$ws.Cells.Item(132, 8).Value = 5287.364
$ws.Cells.Item(132, 9).Value = 3701.9167
$ws.Cells.Item(132, 10).Value = 7189.9
$ws.Cells.Item(132, 11).Value = 11105.7501
$ws.Cells.Item(132, 12).Value = 21569.7
$ws.Cells.Item(132, 13).Value = -8575.750100000001
$ws.Cells.Item(132, 14).Value = -26629.7
